$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '28.569.84'
$r.Style = 'Normal'
$r = $ws.Range('E2')
$r.NumberFormat = '@'
$r.Value = '  +1.73%  '
$r.Style = 'Normal'
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '1.915.06'
$r.Style = 'Normal'
$r = $ws.Range('E3')
$r.NumberFormat = '@'
$r.Value = '  +5.43%  '
$r.Style = 'Normal'
$r = $ws.Range('E4')
$r.NumberFormat = '@'
$r.Value = '  +0.20%  '
$r.Style = 'Normal'
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '315.91'
$r.Style = 'Normal'
$r = $ws.Range('E5')
$r.NumberFormat = '@'
$r.Value = '  +1.77%  '
$r.Style = 'Normal'
$r = $ws.Range('E6')
$r.NumberFormat = '@'
$r.Value = '  +0.16%  '
$r.Style = 'Normal'
$r = $ws.Range('E7')
$r.NumberFormat = '@'
$r.Value = '  +4.58%  '
$r.Style = 'Normal'
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '0.3969'
$r.Style = 'Normal'
$r = $ws.Range('D9')
$r.NumberFormat = '@'
$r.Value = '0.09689'
$r.Style = 'Normal'
$r = $ws.Range('E9')
$r.NumberFormat = '@'
$r.Value = '  -2.19%  '
$r.Style = 'Normal'
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '1.158'
$r.Style = 'Normal'
$r = $ws.Range('E10')
$r.NumberFormat = '@'
$r.Value = '  +4.61%  '
$r.Style = 'Normal'
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '42.13'
$r.Style = 'Normal'
$r = $ws.Range('E11')
$r.NumberFormat = '@'
$r.Value = '  +3.01%  '
$r.Style = 'Normal'
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '6.544'
$r.Style = 'Normal'
$r = $ws.Range('E12')
$r.NumberFormat = '@'
$r.Value = '  +1.88%  '
$r.Style = 'Normal'
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '21.24'
$r.Style = 'Normal'
$r = $ws.Range('E13')
$r.NumberFormat = '@'
$r.Value = '  +3.21%  '
$r.Style = 'Normal'
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '1.919.56'
$r.Style = 'Normal'
$r = $ws.Range('E14')
$r.NumberFormat = '@'
$r.Value = '  +6.03%  '
$r.Style = 'Normal'
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '7.608'
$r.Style = 'Normal'
$r = $ws.Range('E15')
$r.NumberFormat = '@'
$r.Value = '  +4.44%  '
$r.Style = 'Normal'
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '1.001'
$r.Style = 'Normal'
$r = $ws.Range('E16')
$r.NumberFormat = '@'
$r.Value = '  +0.23%  '
$r.Style = 'Normal'
$r = $ws.Range('D17')
$r.NumberFormat = '@'
$r.Value = '0.00001139'
$r.Style = 'Normal'
$r = $ws.Range('E17')
$r.NumberFormat = '@'
$r.Value = '  +0.18%  '
$r.Style = 'Normal'
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '94.22'
$r.Style = 'Normal'
$r = $ws.Range('E18')
$r.NumberFormat = '@'
$r.Value = '  +1.82%  '
$r.Style = 'Normal'
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '0.06660'
$r.Style = 'Normal'
$r = $ws.Range('E19')
$r.NumberFormat = '@'
$r.Value = '  +0.23%  '
$r.Style = 'Normal'
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '18.13'
$r.Style = 'Normal'
$r = $ws.Range('E20')
$r.NumberFormat = '@'
$r.Value = '  +5.39%  '
$r.Style = 'Normal'
$r = $ws.Range('E21')
$r.NumberFormat = '@'
$r.Value = '  +0.23%  '
$r.Style = 'Normal'
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '6.329'
$r.Style = 'Normal'
$r = $ws.Range('E22')
$r.NumberFormat = '@'
$r.Value = '  +6.62%  '
$r.Style = 'Normal'
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '28.648.28'
$r.Style = 'Normal'
$r = $ws.Range('E24')
$r.NumberFormat = '@'
$r.Value = '  +3.11%  '
$r.Style = 'Normal'
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '2.297'
$r.Style = 'Normal'
$r = $ws.Range('E25')
$r.NumberFormat = '@'
$r.Value = '  +1.49%  '
$r.Style = 'Normal'
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '2.691'
$r.Style = 'Normal'
$r = $ws.Range('E26')
$r.NumberFormat = '@'
$r.Value = '  +11.32%  '
$r.Style = 'Normal'
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '2.137.02'
$r.Style = 'Normal'
$r = $ws.Range('E27')
$r.NumberFormat = '@'
$r.Value = '  +5.70%  '
$r.Style = 'Normal'
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '21.28'
$r.Style = 'Normal'
$r = $ws.Range('E28')
$r.NumberFormat = '@'
$r.Value = '  +3.08%  '
$r.Style = 'Normal'
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '159.94'
$r.Style = 'Normal'
$r = $ws.Range('E29')
$r.NumberFormat = '@'
$r.Value = '  +0.89%  '
$r.Style = 'Normal'
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '129.12'
$r.Style = 'Normal'
$r = $ws.Range('E30')
$r.NumberFormat = '@'
$r.Value = '  +1.55%  '
$r.Style = 'Normal'
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '1.109'
$r.Style = 'Normal'
$r = $ws.Range('E31')
$r.NumberFormat = '@'
$r.Value = '  +6.67%  '
$r.Style = 'Normal'
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '0.1086'
$r.Style = 'Normal'
$r = $ws.Range('E32')
$r.NumberFormat = '@'
$r.Value = '  +1.73%  '
$r.Style = 'Normal'
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '5.751'
$r.Style = 'Normal'
$r = $ws.Range('E33')
$r.NumberFormat = '@'
$r.Value = '  +3.13%  '
$r.Style = 'Normal'
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '3.642'
$r.Style = 'Normal'
$r = $ws.Range('E34')
$r.NumberFormat = '@'
$r.Value = '  +1.41%  '
$r.Style = 'Normal'
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '9.914'
$r.Style = 'Normal'
$r = $ws.Range('E35')
$r.NumberFormat = '@'
$r.Value = '  +10.73%  '
$r.Style = 'Normal'
$r = $ws.Range('D36')
$r.NumberFormat = '@'
$r.Value = '0.06778'
$r.Style = 'Normal'
$r = $ws.Range('E36')
$r.NumberFormat = '@'
$r.Value = '  +0.61%  '
$r.Style = 'Normal'
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '0.02432'
$r.Style = 'Normal'
$r = $ws.Range('E37')
$r.NumberFormat = '@'
$r.Value = '  +4.31%  '
$r.Style = 'Normal'
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '1.266'
$r.Style = 'Normal'
$r = $ws.Range('E38')
$r.NumberFormat = '@'
$r.Value = '  +7.26%  '
$r.Style = 'Normal'
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '0.2227'
$r.Style = 'Normal'
$r = $ws.Range('E39')
$r.NumberFormat = '@'
$r.Value = '  +3.86%  '
$r.Style = 'Normal'
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '11.85'
$r.Style = 'Normal'
$r = $ws.Range('E40')
$r.NumberFormat = '@'
$r.Value = '  +4.75%  '
$r.Style = 'Normal'
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '5.117'
$r.Style = 'Normal'
$r = $ws.Range('E41')
$r.NumberFormat = '@'
$r.Value = '  +3.36%  '
$r.Style = 'Normal'
$r = $ws.Range('E42')
$r.NumberFormat = '@'
$r.Value = '  +4.01%  '
$r.Style = 'Normal'
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '1.191'
$r.Style = 'Normal'
$r = $ws.Range('E43')
$r.NumberFormat = '@'
$r.Value = '  +1.36%  '
$r.Style = 'Normal'
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '13.69'
$r.Style = 'Normal'
$r = $ws.Range('E45')
$r.NumberFormat = '@'
$r.Value = '  +3.37%  '
$r.Style = 'Normal'
$r = $ws.Range('D46')
$r.NumberFormat = '@'
$r.Value = '0.6090'
$r.Style = 'Normal'
$r = $ws.Range('E46')
$r.NumberFormat = '@'
$r.Value = '  +2.90%  '
$r.Style = 'Normal'
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '3.778'
$r.Style = 'Normal'
$r = $ws.Range('E47')
$r.NumberFormat = '@'
$r.Value = '  +2.31%  '
$r.Style = 'Normal'
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '1.281'
$r.Style = 'Normal'
$r = $ws.Range('E48')
$r.NumberFormat = '@'
$r.Value = '  -0.08%  '
$r.Style = 'Normal'
$r = $ws.Range('D49')
$r.NumberFormat = '@'
$r.Value = '2.036'
$r.Style = 'Normal'
$r = $ws.Range('E49')
$r.NumberFormat = '@'
$r.Value = '  +5.51%  '
$r.Style = 'Normal'
$r = $ws.Range('D50')
$r.NumberFormat = '@'
$r.Value = '125.43'
$r.Style = 'Normal'
$r = $ws.Range('E50')
$r.NumberFormat = '@'
$r.Value = '  +0.85%  '
$r.Style = 'Normal'
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '1.213'
$r.Style = 'Normal'
$r = $ws.Range('E51')
$r.NumberFormat = '@'
$r.Value = '  +2.62%  '
$r.Style = 'Normal'
